# Adds a new sheet '2025-09-29' at the end of the workbook, containing the
# weekly manga ranking data (header row + 50 ranked rows), mirroring the
# structure of the existing dated sheets.

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Sheets.Count
$lastSheet = $wb.Sheets.Item($sheetCount)
$refSheet = $wb.Sheets.Item(1)
foreach ($s in $wb.Sheets) {
    if ($s.Name -eq "2025-09-22") { $refSheet = $s }
}

$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2025-09-29"

# Header row: rank / title / author / latest_episode, styled like the other sheets.
$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "title"
$ws.Range("C1").Value = "author"
$ws.Range("D1").Value = "latest_episode"
$ws.Range("A1:D1").Style = $refSheet.Range("A1").Style

# Data rows 2..51 (ranks 1..50), written in one shot via a 2D array for speed.
$data = New-Object 'object[,]' 50,4
$data[0,0] = 1
$data[0,1] = '異世界おじさん'
$data[0,2] = '殆ど死んでいる(著者)'
$data[0,3] = '第71話'
$data[1,0] = 2
$data[1,1] = '「おかえり、パパ」'
$data[1,2] = '蝉丸'
$data[1,3] = '第28話　帰宅'
$data[2,0] = 3
$data[2,1] = '悪人面したＢ級冒険者 主人公とその幼馴染たちのパパになる'
$data[2,2] = 'こげめ(著者) えんじ(原作) ハラカズヒロ(キャラクター原案)'
$data[2,3] = '第17話-3：「違法奴隷商討伐」'
$data[3,0] = 4
$data[3,1] = 'ニチアサ好きのオタクが悪役生徒に転生した結果、破滅フラグが崩壊していく件について'
$data[3,2] = '烏丸英（原作） どんぐりす（漫画）'
$data[3,3] = '第14話（後編）急襲…事件の始まり'
$data[4,0] = 5
$data[4,1] = '戸崎さんは僕にだけ冷たい'
$data[4,2] = 'saku(著者)'
$data[4,3] = '第29話-1'
$data[5,0] = 6
$data[5,1] = '魔術師クノンは見えている'
$data[5,2] = 'La-na(作画) 南野海風(原作) Ｌａｒｕｈａ(キャラクター原案)'
$data[5,3] = '第41話①'
$data[6,0] = 7
$data[6,1] = '生徒会にも穴はある！'
$data[6,2] = 'むちまろ'
$data[6,3] = '第136話	ぎゅってしたい'
$data[7,0] = 8
$data[7,1] = 'モブ高生の俺でも冒険者になればリア充になれますか？'
$data[7,2] = '原作：百均 漫画：さぎやまれん キャラクター原案：hai'
$data[7,3] = '第31話'
$data[8,0] = 9
$data[8,1] = '君のラブを見せてくれ！'
$data[8,2] = 'リムコロ(著者)'
$data[8,3] = 'コミックス第⑤巻発売告知'
$data[9,0] = 10
$data[9,1] = '時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―'
$data[9,2] = '光永康則'
$data[9,3] = '第６９話『岩鬼停止』④'
$data[10,0] = 11
$data[10,1] = 'ルパン三世 異世界の姫君（ネイバーワールドプリンセス）'
$data[10,2] = 'モンキー・パンチ／エム・ピー・ワークス 内々けやき 佐伯庸介 白狼'
$data[10,3] = '第111話：泥棒の音を取り戻せ'
$data[11,0] = 12
$data[11,1] = 'クラスで２番目に可愛い女の子と友だちになった'
$data[11,2] = '尾野凛(漫画) たかた(原作) 日向あずり(キャラクター原案)'
$data[11,3] = '第35話②'
$data[12,0] = 13
$data[12,1] = '実は俺、最強でした？'
$data[12,2] = '原作：澄守 彩 漫画：高橋 愛'
$data[12,3] = '第127話　帝国の思惑'
$data[13,0] = 14
$data[13,1] = '男女比1：5の世界でも普通に生きられると思った？　～激重感情な彼女たちが無自覚男子に翻弄されたら～'
$data[13,2] = '三藤 孝太郎(原作) 桃季憂(漫画) jimmy(キャラクター原案)'
$data[13,3] = '第11話-2'
$data[14,0] = 15
$data[14,1] = 'いとこのこ'
$data[14,2] = 'いぬちく(著者)'
$data[14,3] = '休載イラスト'
$data[15,0] = 16
$data[15,1] = '勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～'
$data[15,2] = '漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり'
$data[15,3] = '第５３話　虎の尾を踏む器用貧乏（１）'
$data[16,0] = 17
$data[16,1] = 'ダンジョンの幼なじみ'
$data[16,2] = '久真やすひさ(著者)'
$data[16,3] = '第58話'
$data[17,0] = 18
$data[17,1] = '異世界のんびり農家'
$data[17,2] = '剣康之(作画) 内藤騎之介(原作) やすも(キャラクター原案)'
$data[17,3] = '第306話'
$data[18,0] = 19
$data[18,1] = '【パクパクですわ】追放されたお嬢様の『モンスターを食べるほど強くなる』スキルは、１食で１レベルアップする前代未聞の最強スキルでした。３日で人類最強になりましたわ～！'
$data[18,2] = '島知宏 音速炒飯 有都あらゆる'
$data[18,3] = '第２４食　サーモンのカルパッチョ、パクパクですわ！（２）'
$data[19,0] = 20
$data[19,1] = '美人女上司滝沢さん'
$data[19,2] = 'やんBARU(著者)'
$data[19,3] = '第204話'
$data[20,0] = 21
$data[20,1] = '世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜'
$data[20,2] = '戸賀 環 坂木持丸 riritto'
$data[20,3] = '第53話①　ダンジョンをクリアしてみた'
$data[21,0] = 22
$data[21,1] = 'リビルドワールド'
$data[21,2] = '綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)'
$data[21,3] = '第73話③'
$data[22,0] = 23
$data[22,1] = 'アザミヤコを好きになる'
$data[22,2] = 'ユニティコング(原作) ツノニガウ(作画)'
$data[22,3] = '第10話'
$data[23,0] = 24
$data[23,1] = '男嫌いな美人姉妹を名前も告げずに助けたら一体どうなる?'
$data[23,2] = 'みょん(原作) 司馬淳子(漫画) ぎうにう(キャラクターデザイン)'
$data[23,3] = '第25話'
$data[24,0] = 25
$data[24,1] = '「ククク……。奴は四天王の中でも最弱」と解雇された俺、なぜか勇者と聖女の師匠になる'
$data[24,2] = '漫画：芳橋アツシ 原作：延野正行 キャラクター原案：坂野杏梨'
$data[24,3] = '第42話　奴は座っている（前編）'
$data[25,0] = 26
$data[25,1] = '独身貴族は異世界を謳歌する ～結婚しない男の優雅なおひとりさまライフ～'
$data[25,2] = '漫画：駒鳥 ひわ 原作：錬金王 キャラクター原案：三登 いつき'
$data[25,3] = '第34話 独身貴族は礼の品を贈る（2）'
$data[26,0] = 27
$data[26,1] = '魔導具師ダリヤはうつむかない ～Dahliya Wilts No More～'
$data[26,2] = '漫画：住川惠 原作：甘岸久弥(｢魔導具師ダリヤはうつむかない ～今日から自由な職人ライフ～｣MFブックス刊) キャラクター原案：景、駒田ハチ'
$data[26,3] = '第48話 緑の塔夏祭り夕食会①'
$data[27,0] = 28
$data[27,1] = '物語の黒幕に転生して'
$data[27,2] = '瀬川はじめ(漫画) 結城涼(原作) なかむら(キャラクター原案)'
$data[27,3] = '第34話'
$data[28,0] = 29
$data[28,1] = '姫様“拷問”の時間です'
$data[28,2] = '原作:春原ロビンソン　漫画:ひらけい'
$data[28,3] = '拷問150'
$data[29,0] = 30
$data[29,1] = 'リアデイルの大地にて'
$data[29,2] = '月見だしお(著者) Ceez(原作) てんまそ(キャラクター原案) 涼風涼(構成)'
$data[29,3] = '第40章-2'
$data[30,0] = 31
$data[30,1] = 'おとなりのダウナーさんは無理させない'
$data[30,2] = '瑠璃いろ(著者)'
$data[30,3] = '第13.5話'
$data[31,0] = 32
$data[31,1] = '解雇された暗黒兵士(30代)のスローなセカンドライフ'
$data[31,2] = '岡沢六十四 るれくちぇ sage・ジョー'
$data[31,3] = '第73話(前編) ぶらり新生ラクス街'
$data[32,0] = 33
$data[32,1] = '落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～'
$data[32,2] = '村上よしゆき 茨木野 あるてら'
$data[32,3] = '第４２話　勇者、六邪神将相手に舐めプしてたら、ピンチになる（４）'
$data[33,0] = 34
$data[33,1] = 'Lv２からチートだった元勇者候補のまったり異世界ライフ'
$data[33,2] = '糸町秋音（漫画） 鬼ノ城ミヤ（原作） 片桐（キャラクター原案）'
$data[33,3] = '第61話　成長…？'
$data[34,0] = 35
$data[34,1] = 'アイドル辞めるけど結婚してくれますか!?'
$data[34,2] = '三吉汐美(著者)'
$data[34,3] = '第17話後半'
$data[35,0] = 36
$data[35,1] = 'くらいあの子としたいこと'
$data[35,2] = '碇マナツ(著者)'
$data[35,3] = '第85話'
$data[36,0] = 37
$data[36,1] = '三枝さんはメガネ先輩と恋を描く'
$data[36,2] = 'セレビィ量産型(著者)'
$data[36,3] = '第22話後半'
$data[37,0] = 38
$data[37,1] = '治癒魔法の間違った使い方 ~戦場を駆ける回復要員~'
$data[37,2] = '九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)'
$data[37,3] = '第83話(前編)その1'
$data[38,0] = 39
$data[38,1] = 'ゲーム世界で魔物に転生してしまった俺、前世で推しだったヒロインを拾ってしまう'
$data[38,2] = '三部べべ(漫画) ねうしとら(原作)'
$data[38,3] = '第2話-2'
$data[39,0] = 40
$data[39,1] = 'バーサス'
$data[39,2] = '原作：ONE 漫画：あずま京太郎 構成：bose'
$data[39,3] = '第28話　工場（2）'
$data[40,0] = 41
$data[40,1] = 'インフィニット・デンドログラム'
$data[40,2] = '今井神 原作：海道左近 キャラクター原案：タイキ'
$data[40,3] = '第73話'
$data[41,0] = 42
$data[41,1] = '配信に致命的に向いていない女の子が迷宮で黙々と人助けする配信'
$data[41,2] = '下田将也(漫画) 佐藤悪糖(原作) 福きつね(キャラクター原案)'
$data[41,3] = '第3話前編'
$data[42,0] = 43
$data[42,1] = '老後に備えて異世界で８万枚の金貨を貯めます'
$data[42,2] = 'FUNA 東西 モトエ恵介'
$data[42,3] = '第122話　襲撃［その１］'
$data[43,0] = 44
$data[43,1] = 'お気楽領主の楽しい領地防衛 ～生産系魔術で名もなき村を最強の城塞都市に～'
$data[43,2] = '青色まろ（漫画） 赤池宗（原作） 転（原作イラスト）'
$data[43,3] = '第34話　永住権'
$data[44,0] = 45
$data[44,1] = '幼女戦記'
$data[44,2] = '東條チカ(漫画) カルロ・ゼン(原作) 篠月しのぶ(キャラクター原案)'
$data[44,3] = '第百八章：ドードーバード航空戦Ⅲ'
$data[45,0] = 46
$data[45,1] = '十年目、帰還を諦めた転移者はいまさら主人公になる'
$data[45,2] = '原作：氷純（「十年目、帰還を諦めた転移者はいまさら主人公になる」MFブックス刊） 漫画：しゅーかま キャラクター原案：あんべよしろう'
$data[45,3] = '第１９話②'
$data[46,0] = 47
$data[46,1] = '最果てのパラディン'
$data[46,2] = '奥橋睦（漫画） 柳野かなた（原作） 輪くすさが（キャラクター原案）'
$data[46,3] = '第67話　月の旅路'
$data[47,0] = 48
$data[47,1] = '理想のヒモ生活'
$data[47,2] = '日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)'
$data[47,3] = '第87話　その4'
$data[48,0] = 49
$data[48,1] = 'サーシャちゃんとクラスメイトオタクくん'
$data[48,2] = 'はぐはぐ(著者)'
$data[48,3] = '第86話'
$data[49,0] = 50
$data[49,1] = '異世界でも無難に生きたい症候群'
$data[49,2] = '原作：安泰（一二三書房刊） 漫画：笹峰コウ キャラクター原案：ひたきゆう'
$data[49,3] = '第31話④'

$ws.Range("A2:D51").Value = $data
# Data cells keep the workbook's default (unstyled) cell format, matching
# the other dated sheets' data rows -- no explicit Style assignment needed.

Write-Output "Added sheet '2025-09-29' with $($ws.UsedRange.Rows.Count) rows"
